$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '51.902.38'
$ws.Range("E2").Value = '  -0.47%  '
$ws.Range("D3").Value = '2.785.04'
$ws.Range("E3").Value = '  -1.93%  '
$ws.Range("E4").Value = '  +0.06%  '
$ws.Range("D5").Value = "'357.66"
$ws.Range("E5").Value = '  -1.21%  '
$ws.Range("D6").Value = "'108.84"
$ws.Range("E6").Value = '  -3.69%  '
$ws.Range("D7").Value = "'0.557"
$ws.Range("E7").Value = '  -2.79%  '
$ws.Range("E8").Value = '  +0.02%  '
$ws.Range("D9").Value = "'0.588"
$ws.Range("E9").Value = '  -2.57%  '
$ws.Range("D10").Value = "'40.13"
$ws.Range("E10").Value = '  -3.32%  '
$ws.Range("D11").Value = "'0.0847"
$ws.Range("E11").Value = '  -1.83%  '
$ws.Range("D12").Value = "'0.133"
$ws.Range("E12").Value = '  +1.28%  '
$ws.Range("D13").Value = "'19.43"
$ws.Range("E13").Value = '  -3.01%  '
$ws.Range("D14").Value = "'7.55"
$ws.Range("E14").Value = '  -3.32%  '
$ws.Range("D15").Value = '3.230.85'
$ws.Range("E15").Value = '  -1.61%  '
$ws.Range("D16").Value = '2.780.13'
$ws.Range("E16").Value = '  -1.79%  '
$ws.Range("D17").Value = "'0.939"
$ws.Range("E17").Value = '  +3.24%  '
$ws.Range("D18").Value = '51.872.41'
$ws.Range("E18").Value = '  -0.36%  '
$ws.Range("D19").Value = "'7.51"
$ws.Range("E19").Value = '  -0.88%  '
$ws.Range("D20").Value = "'3.08"
$ws.Range("E20").Value = '  -2.38%  '
$ws.Range("D21").Value = "'13.11"
$ws.Range("E21").Value = '  -3.25%  '
$ws.Range("D22").Value = '0.0₃0974'
$ws.Range("E22").Value = '  -2.38%  '
$ws.Range("D23").Value = "'70.08"
$ws.Range("E23").Value = '  -0.37%  '
$ws.Range("D24").Value = "'269.45"
$ws.Range("E24").Value = '  +0.34%  '
$ws.Range("D25").Value = "'2.75"
$ws.Range("E25").Value = '  -3.15%  '
$ws.Range("D26").Value = "'26.46"
$ws.Range("E26").Value = '  -2.33%  '
$ws.Range("E27").Value = '  -0.08%  '
$ws.Range("E28").Value = '  +16.15%  '
$ws.Range("D29").Value = "'10.27"
$ws.Range("E29").Value = '  -1.44%  '
$ws.Range("D30").Value = "'2.29"
$ws.Range("E30").Value = '  +1.68%  '
$ws.Range("D31").Value = "'0.0471"
$ws.Range("E31").Value = '  -3.49%  '
$ws.Range("D32").Value = "'52.02"
$ws.Range("E32").Value = '  -3.66%  '
$ws.Range("D33").Value = "'34.19"
$ws.Range("E33").Value = '  -3.16%  '
$ws.Range("D34").Value = "'5.72"
$ws.Range("E34").Value = '  -2.58%  '
$ws.Range("D35").Value = "'0.0844"
$ws.Range("E35").Value = '  -0.18%  '
$ws.Range("D36").Value = "'5.19"
$ws.Range("E36").Value = '  -5.63%  '
$ws.Range("E37").Value = '  +0.04%  '
$ws.Range("D38").Value = "'18.81"
$ws.Range("E38").Value = '  +1.98%  '
$ws.Range("D39").Value = "'3.21"
$ws.Range("E39").Value = '  -1.71%  '
$ws.Range("D40").Value = "'1.98"
$ws.Range("E40").Value = '  -4.39%  '
$ws.Range("D41").Value = "'2.63"
$ws.Range("E41").Value = '  +3.84%  '
$ws.Range("D42").Value = "'0.114"
$ws.Range("E42").Value = '  -2.44%  '
$ws.Range("E43").Value = '  -1.18%  '
$ws.Range("D44").Value = "'119.43"
$ws.Range("E44").Value = '  -6.13%  '
$ws.Range("D45").Value = "'21.87"
$ws.Range("E45").Value = '  -7.58%  '
$ws.Range("D46").Value = '2.091.65'
$ws.Range("E46").Value = '  -1.04%  '
$ws.Range("D47").Value = "'3.25"
$ws.Range("E47").Value = '  -5.18%  '
$ws.Range("D49").Value = "'5.72"
$ws.Range("E49").Value = '  -2.95%  '
$ws.Range("D50").Value = "'0.948"
$ws.Range("E50").Value = '  -4.94%  '
$ws.Range("D51").Value = "'8.82"
$ws.Range("E51").Value = '  -3.34%  '
